# Fruta / hortaliza, semanal
# Weekly refresh of price records for rows 18-39 (Maracuyá, Vega Modelo de Temuco):
# dates, volumes, prices and origin are rolled forward with newly reported figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = 44645
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = 30000
$ws.Range("O18").Value = 30000
$ws.Range("P18").Value = 30000
$ws.Range("S18").Value = 1667
$ws.Range("D19").Value = 44438
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("S19").Value = 1944
$ws.Range("D20").Value = 44704
$ws.Range("D21").Value = 44740
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 34000
$ws.Range("O21").Value = 34000
$ws.Range("P21").Value = 34000
$ws.Range("S21").Value = 1889
$ws.Range("D22").Value = 44434
$ws.Range("N22").Value = 35000
$ws.Range("O22").Value = 35000
$ws.Range("P22").Value = 35000
$ws.Range("S22").Value = 1944
$ws.Range("D23").Value = 44448
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 38000
$ws.Range("O23").Value = 38000
$ws.Range("P23").Value = 38000
$ws.Range("S23").Value = 2111
$ws.Range("D24").Value = 44279
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 35000
$ws.Range("O24").Value = 36000
$ws.Range("P24").Value = 35667
$ws.Range("S24").Value = 1982
$ws.Range("D25").Value = 44679
$ws.Range("M25").Value = 35
$ws.Range("N25").Value = 34000
$ws.Range("O25").Value = 34000
$ws.Range("P25").Value = 34000
$ws.Range("R25").Value = "Perú"
$ws.Range("S25").Value = 1889
$ws.Range("M26").Value = 55
$ws.Range("N26").Value = 28000
$ws.Range("O26").Value = 28000
$ws.Range("P26").Value = 28000
$ws.Range("R26").Value = "Región de Arica y Parinacota"
$ws.Range("S26").Value = 1556
$ws.Range("D27").Value = 44377
$ws.Range("M27").Value = 30
$ws.Range("N27").Value = 40000
$ws.Range("O27").Value = 40000
$ws.Range("P27").Value = 40000
$ws.Range("S27").Value = 2222
$ws.Range("D28").Value = 44432
$ws.Range("M28").Value = 10
$ws.Range("N28").Value = 35000
$ws.Range("O28").Value = 35000
$ws.Range("P28").Value = 35000
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 1944
$ws.Range("D29").Value = 44662
$ws.Range("M29").Value = 15
$ws.Range("N29").Value = 30000
$ws.Range("O29").Value = 30000
$ws.Range("P29").Value = 30000
$ws.Range("R29").Value = "Región de Arica y Parinacota"
$ws.Range("S29").Value = 1667
$ws.Range("D30").Value = 44449
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 38000
$ws.Range("O30").Value = 38000
$ws.Range("P30").Value = 38000
$ws.Range("S30").Value = 2111
$ws.Range("D31").Value = 44379
$ws.Range("M31").Value = 10
$ws.Range("N31").Value = 30000
$ws.Range("O31").Value = 30000
$ws.Range("P31").Value = 30000
$ws.Range("S31").Value = 1667
$ws.Range("D32").Value = 44720
$ws.Range("M32").Value = 25
$ws.Range("N32").Value = 34000
$ws.Range("O32").Value = 34000
$ws.Range("P32").Value = 34000
$ws.Range("R32").Value = "Perú"
$ws.Range("S32").Value = 1889
$ws.Range("D33").Value = 44294
$ws.Range("M33").Value = 15
$ws.Range("N33").Value = 35000
$ws.Range("O33").Value = 35000
$ws.Range("P33").Value = 35000
$ws.Range("R33").Value = "Región de Arica y Parinacota"
$ws.Range("S33").Value = 1944
$ws.Range("D34").Value = 44435
$ws.Range("M34").Value = 10
$ws.Range("R34").Value = "Perú"
$ws.Range("M35").Value = 105
$ws.Range("R35").Value = "Región de Arica y Parinacota"
$ws.Range("D36").Value = 44357
$ws.Range("M36").Value = 10
$ws.Range("N36").Value = 38000
$ws.Range("O36").Value = 38000
$ws.Range("P36").Value = 38000
$ws.Range("R36").Value = "Perú"
$ws.Range("S36").Value = 2111
$ws.Range("D37").Value = 44418
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 35000
$ws.Range("O37").Value = 35000
$ws.Range("P37").Value = 35000
$ws.Range("R37").Value = "Región de Arica y Parinacota"
$ws.Range("S37").Value = 1944
$ws.Range("D38").Value = 44726
$ws.Range("N38").Value = 34000
$ws.Range("O38").Value = 34000
$ws.Range("P38").Value = 34000
$ws.Range("S38").Value = 1889
$ws.Range("D39").Value = 44748
